# Add files via upload
# Adds "Matches", "Total 4s", "Total 6s", "Total 50s" and "Total 100s" columns
# (I:M) with per-season stats to the "Seasonal Fact" sheet, then leaves that
# sheet as the active / selected one (it was "Team Dim" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seasonal Fact")

# ---- Header row ----------------------------------------------------------
# I1 keeps the existing bold "header" style used by A1:H1; J1:M1 stay
# unstyled, matching the rest of the new columns.
$ws.Range("I1").Value = "Matches"
$ws.Range("J1").Value = "Total 4s"
$ws.Range("K1").Value = "Total 6s"
$ws.Range("L1").Value = "Total 50s"
$ws.Range("M1").Value = "Total 100s"

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Data rows -------------------------------------------------------------
# Year, Matches, Total 4s, Total 6s, Total 50s, Total 100s
$rows = @(
    @(2008, 58, 1702,  622,  96,  6),
    @(2009, 57, 1316,  506,  82,  2),
    @(2010, 60, 1708,  585,  92,  4),
    @(2011, 73, 1913,  639, 104,  6),
    @(2012, 74, 1911,  733, 118,  6),
    @(2013, 76, 2052,  675, 114,  4),
    @(2014, 60, 1563,  715,  87,  3),
    @(2015, 59, 1607,  692,  89,  4),
    @(2016, 60, 1633,  638,  94,  7),
    @(2017, 59, 1612,  705,  95,  5),
    @(2018, 60, 1674,  869,  98,  5),
    @(2019, 60, 1685,  784,  91,  6),
    @(2020, 60, 1514,  734,  92,  5),
    @(2021, 60, 1452,  687,  88,  4),
    @(2022, 74, 2017, 1062, 112,  8),
    @(2023, 74, 2174, 1124, 153, 12),
    @(2024, 71, 2135, 1260, 148, 14),
    @(2025, 74, 2088, 1215, 154,  9)
)

$r = 2
foreach ($row in $rows) {
    $matches  = $row[1]
    $fours    = $row[2]
    $sixes    = $row[3]
    $fifties  = $row[4]
    $hundreds = $row[5]

    $ws.Cells.Item($r, 9).Value  = $matches    # I - Matches
    $ws.Cells.Item($r, 9).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 10).Value = $fours       # J - Total 4s
    $ws.Cells.Item($r, 10).NumberFormat = "#,##0"

    $ws.Cells.Item($r, 11).Value = $sixes       # K - Total 6s
    if ($sixes -ge 1000) {
        $ws.Cells.Item($r, 11).NumberFormat = "#,##0"
    }

    $ws.Cells.Item($r, 12).Value = $fifties     # L - Total 50s
    $ws.Cells.Item($r, 13).Value = $hundreds    # M - Total 100s

    $r = $r + 1
}

# ---- Selections / active sheet --------------------------------------------
# Before the edit "Team Dim" was the active/selected sheet; after the edit
# "Seasonal Fact" becomes active instead.
$wsTeam = $wb.Worksheets.Item("Team Dim")
$wsTeam.Range("O2").Select()

$ws.Activate()
$ws.Range("H21").Select()
